# Scheduled runner update: refresh computed market/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns H:N) across
# several leve worksheets, matching latest pulled price data.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 56.5
$ws.Range("I33").Value = 56.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 56.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 172.5
$ws.Range("N33").ClearContents()
$ws.Range("H129").Value = 1314.6207
$ws.Range("J129").Value = 1367.7778
$ws.Range("L129").Value = 4103.3334
$ws.Range("N129").Value = -14103.3334

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 927.73334
$ws.Range("I45").Value = 889.1111
$ws.Range("K45").Value = 889.1111
$ws.Range("M45").Value = -512.1111
$ws.Range("H81").Value = 168666.67
$ws.Range("I81").Value = 3000
$ws.Range("J81").Value = 500000
$ws.Range("K81").Value = 3000
$ws.Range("L81").Value = 500000
$ws.Range("M81").Value = -2002
$ws.Range("N81").Value = -501996
$ws.Range("H84").Value = 168666.67
$ws.Range("I84").Value = 3000
$ws.Range("J84").Value = 500000
$ws.Range("K84").Value = 9000
$ws.Range("L84").Value = 1500000
$ws.Range("M84").Value = -4008
$ws.Range("N84").Value = -1509984
$ws.Range("H102").Value = 1920
$ws.Range("I102").Value = 1800
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1800
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -178
$ws.Range("N102").Value = -5244

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 2950
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 2950
$ws.Range("M16").Value = -713
$ws.Range("N16").Value = -3524
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 2950
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 2950
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -7290
$ws.Range("H129").ClearContents()
$ws.Range("I129").ClearContents()
$ws.Range("J129").ClearContents()
$ws.Range("K129").ClearContents()
$ws.Range("L129").ClearContents()
$ws.Range("N129").ClearContents()
$ws.Range("H130").ClearContents()
$ws.Range("I130").ClearContents()
$ws.Range("J130").ClearContents()
$ws.Range("K130").ClearContents()
$ws.Range("L130").ClearContents()
$ws.Range("H131").ClearContents()
$ws.Range("I131").ClearContents()
$ws.Range("J131").ClearContents()
$ws.Range("K131").ClearContents()
$ws.Range("L131").ClearContents()
$ws.Range("H132").ClearContents()
$ws.Range("I132").ClearContents()
$ws.Range("J132").ClearContents()
$ws.Range("K132").ClearContents()
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H133").ClearContents()
$ws.Range("I133").ClearContents()
$ws.Range("J133").ClearContents()
$ws.Range("K133").ClearContents()
$ws.Range("L133").ClearContents()
$ws.Range("M133").ClearContents()
$ws.Range("N133").ClearContents()
$ws.Range("H134").ClearContents()
$ws.Range("I134").ClearContents()
$ws.Range("J134").ClearContents()
$ws.Range("K134").ClearContents()
$ws.Range("L134").ClearContents()
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()
$ws.Range("H135").ClearContents()
$ws.Range("I135").ClearContents()
$ws.Range("J135").ClearContents()
$ws.Range("K135").ClearContents()
$ws.Range("L135").ClearContents()
$ws.Range("N135").ClearContents()
$ws.Range("H137").ClearContents()
$ws.Range("I137").ClearContents()
$ws.Range("J137").ClearContents()
$ws.Range("K137").ClearContents()
$ws.Range("L137").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("H138").ClearContents()
$ws.Range("I138").ClearContents()
$ws.Range("J138").ClearContents()
$ws.Range("K138").ClearContents()
$ws.Range("L138").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("H139").ClearContents()
$ws.Range("I139").ClearContents()
$ws.Range("J139").ClearContents()
$ws.Range("K139").ClearContents()
$ws.Range("L139").ClearContents()
$ws.Range("M139").ClearContents()
$ws.Range("N139").ClearContents()
$ws.Range("H140").ClearContents()
$ws.Range("I140").ClearContents()
$ws.Range("J140").ClearContents()
$ws.Range("K140").ClearContents()
$ws.Range("L140").ClearContents()
$ws.Range("H141").ClearContents()
$ws.Range("I141").ClearContents()
$ws.Range("J141").ClearContents()
$ws.Range("K141").ClearContents()
$ws.Range("L141").ClearContents()
$ws.Range("M141").ClearContents()
$ws.Range("N141").ClearContents()

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 923.26666
$ws.Range("I23").Value = 2570.25
$ws.Range("J23").Value = 324.36365
$ws.Range("K23").Value = 7710.75
$ws.Range("L23").Value = 973.09095
$ws.Range("M23").Value = -7475.75
$ws.Range("N23").Value = -1443.09095
$ws.Range("H75").Value = 2143.4
$ws.Range("I75").Value = 1753
$ws.Range("K75").Value = 5259
$ws.Range("M75").Value = -4261
$ws.Range("H78").Value = 2143.4
$ws.Range("I78").Value = 1753
$ws.Range("K78").Value = 15777
$ws.Range("M78").Value = -10785
$ws.Range("H131").Value = 1446.0492
$ws.Range("I131").Value = 417.77777
$ws.Range("J131").Value = 1624.0193
$ws.Range("K131").Value = 1253.33331
$ws.Range("L131").Value = 4872.0579
$ws.Range("M131").Value = 3786.66669
$ws.Range("N131").Value = -14952.0579

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3142.3157
$ws.Range("I7").Value = 2040.8
$ws.Range("K7").Value = 2040.8
$ws.Range("M7").Value = -1928.8
$ws.Range("H61").Value = 5701.9033
$ws.Range("I61").Value = 5633.0435
$ws.Range("K61").Value = 5633.0435
$ws.Range("M61").Value = -5431.0435
$ws.Range("H113").Value = 5701.9033
$ws.Range("I113").Value = 5633.0435
$ws.Range("K113").Value = 5633.0435
$ws.Range("M113").Value = -3463.0435
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H126").Value = 3142.3157
$ws.Range("I126").Value = 2040.8
$ws.Range("K126").Value = 6122.4
$ws.Range("M126").Value = -3652.4
$ws.Range("H127").Value = 50500
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 50500
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 50500
$ws.Range("N127").Value = -60420
$ws.Range("H128").Value = 41666.668
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 41666.668
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 41666.668
$ws.Range("N128").Value = -51626.668
$ws.Range("H129").Value = 50000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 50000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 50000
$ws.Range("N129").Value = -60000
$ws.Range("H130").Value = 10429
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 10429
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 10429
$ws.Range("N130").Value = -20469
$ws.Range("H131").Value = 50000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 50000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080
$ws.Range("H132").Value = 4326.396
$ws.Range("I132").Value = 3654.5173
$ws.Range("J132").Value = 5351.8945
$ws.Range("K132").Value = 10963.5519
$ws.Range("L132").Value = 16055.6835
$ws.Range("M132").Value = -8433.5519
$ws.Range("N132").Value = -21115.6835
$ws.Range("H133").Value = 50811.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 50811.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 50811.5
$ws.Range("N133").Value = -55871.5
$ws.Range("H134").Value = 7000
$ws.Range("I134").Value = 7000
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7000
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1930
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H136").Value = 4052.2273
$ws.Range("I136").Value = 2362.1516
$ws.Range("J136").Value = 9122.454
$ws.Range("K136").Value = 7086.4548
$ws.Range("L136").Value = 27367.362
$ws.Range("M136").Value = -4536.4548
$ws.Range("N136").Value = -32467.362
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 44828.332
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 44828.332
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 44828.332
$ws.Range("N139").Value = -55108.332
$ws.Range("H140").Value = 70000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 70000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 70000
$ws.Range("N140").Value = -80360
$ws.Range("H141").Value = 55914.285
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 55914.285
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 55914.285
$ws.Range("N141").Value = -66274.285
